$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "cleiton.souza@mrv.com.br"
$ws.Range("B6").Value = "Planilha de medição de MOP - v02"
$ws.Range("C6").Value = "Ferramenta de Planejamento"
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = "Teste 1 - usuário 2`n"
